$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value / type changes -------------------------------------------
# "Position" column values were stored as text ("1"/"2"); make them real numbers.
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C6").Value = 1

# Fill in the previously-empty multi-line description for the second row of
# the first file group.
$ws.Range("G4").Value = "some description" + [char]10 + "with multiple lines"

# --- Column widths ---------------------------------------------------------
# Switch from auto-fit widths to fixed widths (no more bestFit).
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 13.166666666666666
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

# --- Wrap text --------------------------------------------------------
# All the used styles should wrap long text (needed for the new
# multi-line description cell).
$ws.Range("A1:H6").WrapText = $true

# Keep the row at its natural (non-custom) height even though row 4 now
# holds a 2-line value - match the original row-height metadata.
$ws.Rows.Item(4).AutoFit()

# --- Fill / header color -----------------------------------------------
# Darken the header-row highlight fill color (RGB 88,88,88 == 0xFF585858).
$ws.Range("A1:H1").Interior.Color = 5789784

# --- AutoFilter ----------------------------------------------------------
$ws.Range("A1:H7").AutoFilter()
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet0!`$A`$1:`$H`$7")
$fd.Visible = $false
